# SwaadSutra Daily Orders workbook update
# New order placed: Pooja / a14 / Wheat Chapati x1 / Total 15 / NEW / PENDING
# logged at 2026-01-13 16:54, pushed onto the top of the "Daily Orders" log
# (order id 7), bumping every later row down by one.

$wb = $excel.ActiveWorkbook

# ---- 1. "Daily Orders" sheet -------------------------------------------
$ws = $wb.Worksheets.Item("Daily Orders")

# Insert a fresh row 2 - this shifts the existing order history down
# (row 2->3, 3->4, ... 7->8) and keeps their values/formatting intact.
$ws.Rows(2).Insert()

# Fill in the details of the brand-new order in the now-empty row 2.
$ws.Range("A2").Value = 7
$ws.Range("B2").Value = "2026-01-13 16:54"
$ws.Range("C2").Value = "Pooja"
$ws.Range("D2").Value = "a14"

# Phone number reads as text (matches the rest of the sheet's Phone column).
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "9096648553"

$ws.Range("F2").Value = "Wheat Chapati x1"
$ws.Range("G2").Value = 15
$ws.Range("H2").Value = "NEW"
$ws.Range("I2").Value = "PENDING"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = ""

# ---- 2. "Summary" sheet --------------------------------------------------
$sum = $wb.Worksheets.Item("Summary")
$sum.Range("A2").Value = 7          # Total Orders      6 -> 7
$sum.Range("B2").Value = 6          # New               5 -> 6
$sum.Range("G2").Value = 195        # Total Revenue   180 -> 195

# ---- 3. "Items Breakdown" sheet ------------------------------------------
$items = $wb.Worksheets.Item("Items Breakdown")
$items.Range("B2").Value = 3        # Wheat Chapati qty   2 -> 3
$items.Range("C2").Value = 45       # Wheat Chapati revenue 30 -> 45
